$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.099.60'
$ws.Range('E2').Value = '  +2.36%  '
$ws.Range('D3').Value = '3.062.35'
$ws.Range('E3').Value = '  +1.50%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '559.06'
$ws.Range('E5').Value = '  +2.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.29'
$ws.Range('E6').Value = '  +2.24%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').Value = '3.062.14'
$ws.Range('E8').Value = '  +1.53%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.513'
$ws.Range('E9').Value = '  +4.33%  '
$ws.Range('E10').Value = '  +4.75%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.09'
$ws.Range('E11').Value = '  -12.85%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.483'
$ws.Range('E12').Value = '  +7.52%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000232'
$ws.Range('E13').Value = '  +4.66%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.59'
$ws.Range('E14').Value = '  +3.88%  '
$ws.Range('D15').Value = '3.556.94'
$ws.Range('E15').Value = '  +2.13%  '
$ws.Range('D16').Value = '64.069.58'
$ws.Range('E16').Value = '  +2.44%  '
$ws.Range('D17').Value = '3.061.75'
$ws.Range('E17').Value = '  +1.73%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.77'
$ws.Range('E19').Value = '  +2.51%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '479.42'
$ws.Range('E20').Value = '  +1.13%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.01'
$ws.Range('E21').Value = '  +3.77%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.682'
$ws.Range('E22').Value = '  +3.72%  '
$ws.Range('B23').Value = 'InternetComputer(DFINITY)'
$ws.Range('C23').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '14.49'
$ws.Range('E23').Value = '  +14.10%  '
$ws.Range('B24').Value = 'Uniswap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.63'
$ws.Range('E24').Value = '  +5.73%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '82.16'
$ws.Range('E25').Value = '  +2.96%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.37%  '
$ws.Range('E27').Value = '  +2.38%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.01'
$ws.Range('E28').Value = '  +3.98%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.04'
$ws.Range('E29').Value = '  +0.77%  '
$ws.Range('E30').Value = '  +0.24%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '26.32'
$ws.Range('E31').Value = '  +2.74%  '
$ws.Range('E32').Value = '  +1.01%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.44'
$ws.Range('E33').Value = '  +3.17%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.70'
$ws.Range('E34').Value = '  +1.60%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.24'
$ws.Range('E35').Value = '  +6.12%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '54.90'
$ws.Range('E36').Value = '  +0.36%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0409'
$ws.Range('E37').Value = '  +3.54%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '445.38'
$ws.Range('E38').Value = '  -2.18%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0811'
$ws.Range('E39').Value = '  -0.54%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.87'
$ws.Range('E40').Value = '  +11.37%  '
$ws.Range('D41').Value = '2.991.92'
$ws.Range('E41').Value = '  +0.76%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.27'
$ws.Range('E42').Value = '  +1.94%  '
$ws.Range('E43').Value = '  +0.19%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '28.02'
$ws.Range('E44').Value = '  +3.68%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.262'
$ws.Range('E45').Value = '  +4.71%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.16'
$ws.Range('E46').Value = '  +6.94%  '
$ws.Range('E47').Value = '  -0.02%  '
$ws.Range('E48').Value = '  +3.58%  '
$ws.Range('B49').Value = 'PEPE'
$ws.Range('C49').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D49').Value = '0.0₃0522'
$ws.Range('E49').Value = '  +4.70%  '
$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '118.92'
$ws.Range('E50').Value = '  +2.96%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.08'
$ws.Range('E51').Value = '  +2.43%  '
